$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple +4 bumps to column B (Taxonsorteringsordning)
$ws.Range("B2").Value = 79243
$ws.Range("B3").Value = 92267
$ws.Range("B4").Value = 91808
$ws.Range("B5").Value = 92267
$ws.Range("B6").Value = 92267
$ws.Range("B7").Value = 91771

# Row 8: becomes the former row 9's species record (Tretåig hackspett / Picoides tridactylus)
$ws.Range("A8").Value = 130834387
$ws.Range("B8").Value = 57884
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = "Tretåig hackspett"
$ws.Range("G8").Value = "Picoides tridactylus"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("M8").Value = "äldre spår"
$ws.Range("Q8").Value = 424517
$ws.Range("R8").Value = 6711378

# Row 9: becomes the former row 8's species record (Ullticka / Phellinidium ferrugineofuscum)
$ws.Range("A9").Value = 130834377
$ws.Range("B9").Value = 91808
$ws.Range("E9").Value = 1202
$ws.Range("F9").Value = "Ullticka"
$ws.Range("G9").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H9").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M9").ClearContents()
$ws.Range("Q9").Value = 424489
$ws.Range("R9").Value = 6711391

# Simple +4 bumps to column B
$ws.Range("B10").Value = 91808
$ws.Range("B11").Value = 91808
$ws.Range("B12").Value = 92267
$ws.Range("B13").Value = 92267
$ws.Range("B14").Value = 91808

# Row 15: becomes the former row 16's species record (Vedticka / Fuscoporia viticola)
$ws.Range("A15").Value = 130834374
$ws.Range("B15").Value = 91771
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 5447
$ws.Range("F15").Value = "Vedticka"
$ws.Range("G15").Value = "Fuscoporia viticola"
$ws.Range("H15").Value = "(Schwein.) Murrill"
$ws.Range("Q15").Value = 424489
$ws.Range("R15").Value = 6711391

# Row 16: becomes the former row 15's species record (Ullticka / Phellinidium ferrugineofuscum)
$ws.Range("A16").Value = 130834380
$ws.Range("B16").Value = 91808
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 1202
$ws.Range("F16").Value = "Ullticka"
$ws.Range("G16").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H16").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q16").Value = 424490
$ws.Range("R16").Value = 6711347
